$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '27.972.88'
$ws.Range("E2").Value = '  -0.54%  '

# Row 3
$ws.Range("D3").Value = '1.857.41'
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
Set-TextValue $ws "D4" '1.004'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
Set-TextValue $ws "D5" '312.23'
$ws.Range("E5").Value = '  -0.58%  '

# Row 6
$ws.Range("E6").Value = '  -0.08%  '

# Row 7
Set-TextValue $ws "D7" '0.5129'
$ws.Range("E7").Value = '  +1.44%  '

# Row 8
Set-TextValue $ws "D8" '0.3829'
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
Set-TextValue $ws "D9" '0.08220'
$ws.Range("E9").Value = '  -4.06%  '

# Row 10
Set-TextValue $ws "D10" '1.109'
$ws.Range("E10").Value = '  -0.60%  '

# Row 11
Set-TextValue $ws "D11" '41.47'
$ws.Range("E11").Value = '  -0.14%  '

# Row 12
Set-TextValue $ws "D12" '6.177'
$ws.Range("E12").Value = '  -2.34%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws "D13" '20.50'
$ws.Range("E13").Value = '  -0.74%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.859.06'
$ws.Range("E14").Value = '  -1.88%  '

# Row 15
Set-TextValue $ws "D15" '7.233'
$ws.Range("E15").Value = '  +1.01%  '

# Row 16
Set-TextValue $ws "D16" '1.003'
$ws.Range("E16").Value = '  -0.11%  '

# Row 17
Set-TextValue $ws "D17" '0.00001095'
$ws.Range("E17").Value = '  -0.46%  '

# Row 18
Set-TextValue $ws "D18" '90.38'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
Set-TextValue $ws "D19" '0.06644'
$ws.Range("E19").Value = '  +0.47%  '

# Row 20
$ws.Range("E20").Value = '  -2.77%  '

# Row 21
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
Set-TextValue $ws "D22" '6.006'
$ws.Range("E22").Value = '  -1.44%  '

# Row 23
$ws.Range("D23").Value = '28.007.33'
$ws.Range("E23").Value = '  -0.53%  '

# Row 24
Set-TextValue $ws "D24" '11.04'
$ws.Range("E24").Value = '  -3.10%  '

# Row 25
Set-TextValue $ws "D25" '2.254'
$ws.Range("E25").Value = '  -0.97%  '

# Row 26
$ws.Range("D26").Value = '2.075.41'
$ws.Range("E26").Value = '  -1.57%  '

# Row 27
Set-TextValue $ws "D27" '2.499'
$ws.Range("E27").Value = '  -2.23%  '

# Row 28
Set-TextValue $ws "D28" '157.11'
$ws.Range("E28").Value = '  -0.21%  '

# Row 29
Set-TextValue $ws "D29" '20.40'
$ws.Range("E29").Value = '  -1.53%  '

# Row 30
Set-TextValue $ws "D30" '124.52'
$ws.Range("E30").Value = '  -2.01%  '

# Row 31
Set-TextValue $ws "D31" '0.1065'
$ws.Range("E31").Value = '  +1.52%  '

# Row 32
Set-TextValue $ws "D32" '1.024'
$ws.Range("E32").Value = '  -3.20%  '

# Row 33
Set-TextValue $ws "D33" '5.910'
$ws.Range("E33").Value = '  +5.70%  '

# Row 34
Set-TextValue $ws "D34" '3.595'
$ws.Range("E34").Value = '  -0.29%  '

# Row 35
Set-TextValue $ws "D35" '9.372'
$ws.Range("E35").Value = '  -3.01%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D36" '0.02410'
$ws.Range("E36").Value = '  -1.16%  '

# Row 37
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D37" '0.06490'
$ws.Range("E37").Value = '  -0.86%  '

# Row 38
Set-TextValue $ws "D38" '0.2177'
$ws.Range("E38").Value = '  +0.35%  '

# Row 39
$ws.Range("E39").Value = '  +2.56%  '

# Row 40
Set-TextValue $ws "D40" '1.189'
$ws.Range("E40").Value = '  -1.08%  '

# Row 41
Set-TextValue $ws "D41" '4.968'
$ws.Range("E41").Value = '  +1.41%  '

# Row 42
Set-TextValue $ws "D42" '1.209'
$ws.Range("E42").Value = '  -2.68%  '

# Row 43
Set-TextValue $ws "D43" '11.13'
$ws.Range("E43").Value = '  -3.40%  '

# Row 44
Set-TextValue $ws "D44" '0.6131'
$ws.Range("E44").Value = '  +2.55%  '

# Row 45
Set-TextValue $ws "D45" '12.94'
$ws.Range("E45").Value = '  -1.68%  '

# Row 46
Set-TextValue $ws "D46" '1.277'
$ws.Range("E46").Value = '  -0.02%  '

# Row 47
Set-TextValue $ws "D47" '3.671'
$ws.Range("E47").Value = '  -0.18%  '

# Row 48
Set-TextValue $ws "D48" '2.005'
$ws.Range("E48").Value = '  +1.05%  '

# Row 49
Set-TextValue $ws "D49" '1.212'
$ws.Range("E49").Value = '  -1.54%  '

# Row 50
Set-TextValue $ws "D50" '120.66'
$ws.Range("E50").Value = '  -0.44%  '

# Row 51
Set-TextValue $ws "D51" '78.07'
$ws.Range("E51").Value = '  -2.37%  '
